$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ReferenceFile")
$ws.Range("A3").Value = "smaht:reference_file-B"
$ws.Range("A2").Value = "smaht:reference_file-A"
$ws.Activate()
$ws.Range("A2").Select()
